# Add the new attendance record (row 11) to the sheet, matching the
# existing rows' text formatting (all values - including numeric-looking
# group/year/series codes - are stored as plain text, like "1", "2", "3").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11:I11").NumberFormat = "@"

$ws.Range("A11").Value = "Alexia"
$ws.Range("B11").Value = "ioana-alexia.badea27@s.fpse.unibuc.ro"
$ws.Range("C11").Value = "25.04.2025"
$ws.Range("D11").Value = "21:44:30"
$ws.Range("E11").Value = "1"
$ws.Range("F11").Value = "2"
$ws.Range("G11").Value = "3"
$ws.Range("H11").Value = "TMI II"
$ws.Range("I11").Value = "Seminar"
